$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Audit Reports")

# Row 2 (CC001 / NIIT Delhi) audit record was reset back to "Not Submitted":
#  - Audit Date bumped from 23/11/2025 to 25/11/2025
#  - Audit Data JSON: FO3 remark "hello kais ho aap" cleared to ""
#  - Status / Approval fields cleared/reset
$ws.Range("L2").Value = "25/11/2025"
$ws.Range("M2").Value = '{"FO1":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"FO2":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"FO3":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"FO4":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"FO5":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP1":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP2":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP3":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP4":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP5":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP6":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP7":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP8":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP9":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"DP10":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"PP1":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"PP2":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"MP1":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"MP2":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"MP3":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"MP4":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""},"MP5":{"totalSamples":"","samplesCompliant":"","compliantPercent":0,"score":0,"remarks":""}}'
$ws.Range("N2").Value = "Not Submitted"
$ws.Range("O2").Value = "Not Submitted"
$ws.Range("P2").Value = ""
$ws.Range("Q2").Value = ""
$ws.Range("R2").Value = ""
